# Tkinter Email list + attaching files
# Append 6 new email addresses (with mailto: hyperlinks) below the existing list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$emails = @(
    "d_kirol@lesta.group",
    "d_kireev@lesta.group",
    "d_kinov@lesta.group",
    "d_karin@lesta.group",
    "a_dyshina@lesta.group",
    "e_semenova@lesta.group"
)

$row = 16
foreach ($email in $emails) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $email
    $ws.Hyperlinks.Add($cell, "mailto:$email")
    # Re-apply the Hyperlink cell style by name so it reuses the workbook's
    # existing "Hyperlink" style slot (same one used by the rows above)
    # instead of leaving the cell tagged with the ad-hoc style the Add()
    # call creates on its own.
    $cell.Style = "Hyperlink"
    $row++
}

# Move the active selection to the next empty row, like Excel does after data entry.
$ws.Cells.Item($row, 1).Select() | Out-Null
